$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing S1..S5 step table (rows 10-14) ---
$ws.Range("B11").Value = "S2"
$ws.Range("B12").Value = "S3"
$ws.Range("B13").Value = "S4"
$ws.Range("B14").Value = "S5"

$ws.Range("C10").Value = '''= new Y[] { new Y("a1") } '
$ws.Range("C11").Value = '''= $S1[(X x) transform to x.name]'
$ws.Range("C12").Value = '''= $S1[(X x) transform unique to x.name]'
$ws.Range("C13").Value = '''= $S1[(X x) select all having x.name.length > 0]'
$ws.Range("C14").Value = '''= $S1[(X x) split by x.name.length > 0]'

# --- New "Datatype X" table (rows 21-22) ---
$ws.Range("B21:C21").Merge()
$ws.Range("B21:C21").HorizontalAlignment = -4108
$ws.Range("B21:C21").VerticalAlignment = -4108
$ws.Range("B21:C21").WrapText = $true
$ws.Range("B21").Value = "Datatype X"

$ws.Range("B22:C22").HorizontalAlignment = -4108
$ws.Range("B22:C22").VerticalAlignment = -4108
$ws.Range("B22:C22").WrapText = $true
$ws.Range("B22").Value = "String"
$ws.Range("C22").Value = "name"

# --- New "Datatype Y" table (rows 25-26) ---
$ws.Range("B25:C25").Merge()
$ws.Range("B25:C25").HorizontalAlignment = -4108
$ws.Range("B25:C25").VerticalAlignment = -4108
$ws.Range("B25:C25").WrapText = $true
$ws.Range("B25").Value = "Datatype Y"

$ws.Range("B26:C26").HorizontalAlignment = -4108
$ws.Range("B26:C26").VerticalAlignment = -4108
$ws.Range("B26:C26").WrapText = $true
$ws.Range("B26").Value = "String"
$ws.Range("C26").Value = "name"

# --- Column width ---
$ws.Columns.Item(2).ColumnWidth = 14.666666666666666

# --- Sheet view / selection ---
$ws.Range("F7").Select()
